$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1494059
$ws.Range("J17").Value = 1540009
$ws.Range("L17").Value = 4620027
$ws.Range("N17").Value = -4620363
$ws.Range("H18").Value = 516.3333
$ws.Range("J18").Value = 700
$ws.Range("L18").Value = 700
$ws.Range("N18").Value = -1268
$ws.Range("H40").Value = 1786.3077
$ws.Range("I40").Value = 1133.3334
$ws.Range("J40").Value = 2346
$ws.Range("K40").Value = 1133.3334
$ws.Range("L40").Value = 2346
$ws.Range("M40").Value = -958.3334
$ws.Range("N40").Value = -2696
$ws.Range("H132").Value = 6892.154
$ws.Range("I132").Value = 9448.25
$ws.Range("K132").Value = 28344.75
$ws.Range("M132").Value = -25814.75
$ws.Range("H137").Value = 42464.64
$ws.Range("I137").Value = 2685.3125
$ws.Range("K137").Value = 8055.9375
$ws.Range("M137").Value = -5505.9375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1215.2106
$ws.Range("I99").Value = 839.26666
$ws.Range("K99").Value = 839.26666
$ws.Range("M99").Value = 658.73334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 11
$ws.Range("I25").Value = 11
$ws.Range("K25").Value = 11
$ws.Range("M25").Value = 163
$ws.Range("H31").Value = 13950.767
$ws.Range("I31").Value = 42014
$ws.Range("J31").Value = 3745.9546
$ws.Range("K31").Value = 42014
$ws.Range("L31").Value = 3745.9546
$ws.Range("M31").Value = -41719
$ws.Range("N31").Value = -4335.9546
$ws.Range("H34").Value = 13950.767
$ws.Range("I34").Value = 42014
$ws.Range("J34").Value = 3745.9546
$ws.Range("K34").Value = 42014
$ws.Range("L34").Value = 3745.9546
$ws.Range("M34").Value = -41812
$ws.Range("N34").Value = -4149.9546
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()
$ws.Range("H58").Value = 13188
$ws.Range("I58").Value = 1045.25
$ws.Range("K58").Value = 1045.25
$ws.Range("M58").Value = -842.25
$ws.Range("H60").Value = 11861.789
$ws.Range("I60").Value = 5000
$ws.Range("J60").Value = 12243
$ws.Range("K60").Value = 5000
$ws.Range("L60").Value = 12243
$ws.Range("M60").Value = -4489
$ws.Range("N60").Value = -13265
$ws.Range("H94").Value = 2595.9
$ws.Range("I94").Value = 1609.091
$ws.Range("J94").Value = 3802
$ws.Range("K94").Value = 1609.091
$ws.Range("L94").Value = 3802
$ws.Range("M94").Value = -1158.091
$ws.Range("N94").Value = -4704
$ws.Range("H132").Value = 16505.527
$ws.Range("I132").Value = 20342.963
$ws.Range("J132").Value = 4993.222
$ws.Range("K132").Value = 61028.889
$ws.Range("L132").Value = 14979.666
$ws.Range("M132").Value = -58498.889
$ws.Range("N132").Value = -20039.666
$ws.Range("H134").Value = 1112.8679
$ws.Range("I134").Value = 813.62964
$ws.Range("J134").Value = 1423.6154
$ws.Range("K134").Value = 2440.88892
$ws.Range("L134").Value = 4270.8462
$ws.Range("M134").Value = 94.11108000000013
$ws.Range("N134").Value = -9340.8462
$ws.Range("H136").Value = 13188
$ws.Range("I136").Value = 1045.25
$ws.Range("K136").Value = 3135.75
$ws.Range("M136").Value = -585.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 7945.231
$ws.Range("J107").Value = 273.5
$ws.Range("L107").Value = 820.5
$ws.Range("N107").Value = -4660.5
$ws.Range("H123").Value = 4502.5
$ws.Range("I123").Value = 1510
$ws.Range("J123").Value = 7495
$ws.Range("K123").Value = 4530
$ws.Range("L123").Value = 22485
$ws.Range("M123").Value = -2080
$ws.Range("N123").Value = -27385
$ws.Range("H131").Value = 792.2
$ws.Range("J131").Value = 797.1134
$ws.Range("L131").Value = 2391.3402
$ws.Range("N131").Value = -12471.3402

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 50626.25
$ws.Range("I132").Value = 51847
$ws.Range("J132").Value = 48295.727
$ws.Range("K132").Value = 155541
$ws.Range("L132").Value = 144887.181
$ws.Range("M132").Value = -153011
$ws.Range("N132").Value = -149947.181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 999.8570999999999
$ws.Range("I46").Value = 728.2857
$ws.Range("J46").Value = 1271.4286
$ws.Range("K46").Value = 728.2857
$ws.Range("L46").Value = 1271.4286
$ws.Range("M46").Value = -540.2857
$ws.Range("N46").Value = -1647.4286
$ws.Range("H68").Value = 3055.8
$ws.Range("I68").Value = 3133.3333
$ws.Range("J68").Value = 2939.5
$ws.Range("K68").Value = 3133.3333
$ws.Range("L68").Value = 2939.5
$ws.Range("M68").Value = -2384.3333
$ws.Range("N68").Value = -4437.5
$ws.Range("H71").Value = 3055.8
$ws.Range("I71").Value = 3133.3333
$ws.Range("J71").Value = 2939.5
$ws.Range("K71").Value = 15666.6665
$ws.Range("L71").Value = 14697.5
$ws.Range("M71").Value = -11922.6665
$ws.Range("N71").Value = -22185.5
$ws.Range("H130").Value = 9800
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 9800
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 9800
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -19840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 40000
$ws.Range("J46").Value = 40000
$ws.Range("L46").Value = 40000
$ws.Range("N46").Value = -40462
$ws.Range("H70").Value = 15760
$ws.Range("J70").Value = 15760
$ws.Range("L70").Value = 15760
$ws.Range("N70").Value = -16390
$ws.Range("H73").Value = 15760
$ws.Range("J73").Value = 15760
$ws.Range("L73").Value = 15760
$ws.Range("N73").Value = -17944
$ws.Range("H134").Value = 40000
$ws.Range("J134").Value = 40000
$ws.Range("L134").Value = 120000
$ws.Range("N134").Value = -125070
$ws.Range("H136").Value = 23257392
$ws.Range("I136").Value = 45456144
$ws.Range("J136").Value = 1554
$ws.Range("K136").Value = 136368432
$ws.Range("L136").Value = 4662
$ws.Range("M136").Value = -136365882
$ws.Range("N136").Value = -9762
